$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.540.85'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +2.18%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.673.89'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +1.85%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '220.59'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.63%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5288'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +2.98%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06389'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.86'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +5.45%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07804'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.670.86'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.65%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5568'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.85%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0₅8344'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '65.68'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +1.49%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.534.56'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('E18').Value = '  -0.08%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.767'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +1.36%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '193.53'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.36'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.92%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.315'
$ws.Range('D22').Style = "Normal"
$ws.Range('E24').Value = '  +4.17%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '139.57'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -3.29%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.412'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.33'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.88%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.429'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +2.74%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.06229'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +5.15%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.294'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +2.37%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.620'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +6.54%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.431'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.14%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.685'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.29%  '
$ws.Range('E34').Value = '  +1.65%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.6111'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +8.46%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.413'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.779'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('E38').Value = '  +0.92%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.060'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +3.25%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.090.84'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +5.86%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8598'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.88%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '100.72'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.92%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.818.68'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.46%  '
$ws.Range('E45').Value = '  +4.44%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '58.48'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +5.11%  '
$ws.Range('E47').Value = '  +0.32%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.132'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.523'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +11.09%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05198'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +1.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.017'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +2.11%  '
